$d = $word.ActiveDocument

# --- Step 1: split "Test 1" into two runs "Test " and "1" within the
# same paragraph (no formatting difference between them, just a run
# boundary). This is achieved by inserting a paragraph break between
# "Test " and "1" (which forces a run split at that point) and then
# immediately removing the paragraph mark again, which re-joins the two
# paragraphs back into one while leaving the run split intact.
$findRng = $d.Content
$findRng.Find.ClearFormatting()
[void]$findRng.Find.Execute("Test 1", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$splitPos = $findRng.Start + 5   # position between "Test " and "1"

$splitRng = $d.Range($splitPos, $splitPos)
$splitRng.InsertParagraphAfter()

$para1 = $d.Paragraphs(1).Range
$markRng = $d.Range($para1.End - 1, $para1.End)
$markRng.Delete()

# --- Step 2: the "_GoBack" bookmark currently sits right after "1" (end
# of paragraph 1). Remove it for now; it will be re-created at the end of
# the new "Test 2" paragraph we are about to add.
$goBack = $d.Bookmarks("_GoBack")
$goBack.Delete()

# --- Step 3: append a new, empty paragraph at the end of the document.
$endPos = $d.Content.End
$endRng = $d.Range($endPos - 1, $endPos - 1)
$endRng.InsertParagraphAfter()

# --- Step 4: fill the new paragraph with "Test 2X" (note the temporary
# trailing "X"). Adding the bookmark exactly at the end of a paragraph's
# text (i.e. a collapsed range abutting the paragraph mark) snaps it to
# wrap the wrong paragraph in this runtime, so we keep a placeholder
# character after the insertion point while we create the bookmark, then
# delete the placeholder afterwards.
$newParaCount = $d.Paragraphs.Count
$newPara = $d.Paragraphs($newParaCount).Range
$newPara.InsertBefore("Test 2X")

$newPara = $d.Paragraphs($newParaCount).Range
$bmPos = $newPara.Start + 6   # right after "Test 2", before the "X"
$bmRng = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRng)

# --- Step 5: remove the temporary placeholder character.
$placeholder = $d.Range($bmPos, $bmPos + 1)
$placeholder.Delete()
